#
# Domain Model Class details Stock_details_spec
# Adds a "Packages" column (F) mapping each Class/Service to its .NET
# namespace, renames the "InventoriesService" class to "ProductsService",
# adds hyperlinks on D2:D4, and styles the header row (C1:F1) bold with a
# themed fill (and a bottom border on C1:E1).
#

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rewrite the D / E columns (rows 2-15) with the final class / layer
#    values (row 4 changes from InventoriesService -> ProductsService,
#    and every row below it shifts up by one vs. the original sheet).
# ---------------------------------------------------------------------

$rows = @(
    @{ Row = 2;  D = "StockService";                                       E = "Web Services" },
    @{ Row = 3;  D = "ShipmentsService";                                   E = "Domain Services" },
    @{ Row = 4;  D = "ProductsService";                                    E = "Domain Services" },
    @{ Row = 5;  D = "ShipmentReceipt";                                    E = "Domain Models" },
    @{ Row = 6;  D = "PackagingContent";                                   E = "Domain Models" },
    @{ Row = 7;  D = "InventoryItem";                                      E = "Domain Models" },
    @{ Row = 8;  D = "InventorytItemShipmentReceipt";                      E = "Domain Models" },
    @{ Row = 9;  D = "InventorytItemShipmentReceiptStatus";                E = "Domain Models" },
    @{ Row = 10; D = "ShipmentPackageRepository";                          E = "Domain Infrastructures" },
    @{ Row = 11; D = "PackagingContentRepository";                         E = "Domain Infrastructures" },
    @{ Row = 12; D = "InventoryItemRepository";                            E = "Domain Infrastructures" },
    @{ Row = 13; D = "InventorytItemShipmentReceiptStatusTypeRepository";  E = "Domain Infrastructures" },
    @{ Row = 14; D = "InventorytItemShipmentReceiptStatusRepository";      E = "Domain Infrastructures" },
    @{ Row = 15; D = "InventorytItemShipmentReceiptRepository";            E = "Domain Infrastructures" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# ---------------------------------------------------------------------
# 2. New column F ("Packages") - the .NET namespace each row belongs to.
# ---------------------------------------------------------------------

$ws.Range("F1").Value = "Packages"

$packages = @(
    @{ Row = 2;  F = "Sahapat.StockService.Interfaces" },
    @{ Row = 3;  F = "Sahapat.Shipments.IServices" },
    @{ Row = 4;  F = "Sahapat.Products.IServices" },
    @{ Row = 5;  F = "Sahapat.Shipments.Models" },
    @{ Row = 6;  F = "Sahapat.Shipments.Models" },
    @{ Row = 7;  F = "Sahapat.Products.Models" },
    @{ Row = 8;  F = "Sahapat.Products.Models" },
    @{ Row = 9;  F = "Sahapat.Products.Models" },
    @{ Row = 10; F = "Sahapat.Shipments.IRepositories" },
    @{ Row = 11; F = "Sahapat.Shipments.IRepositories" },
    @{ Row = 12; F = "Sahapat.Products.IRepositories" },
    @{ Row = 13; F = "Sahapat.Products.IRepositories" },
    @{ Row = 14; F = "Sahapat.Products.IRepositories" },
    @{ Row = 15; F = "Sahapat.Products.IRepositories" }
)

foreach ($r in $packages) {
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}

$ws.Columns.Item(6).ColumnWidth = 35.16666666666667

# ---------------------------------------------------------------------
# 3. Hyperlinks for D2 / D3 / D4 (draw.io class-diagram deep links).
#    D4's cached display text stays "InventoriesService" (stale from the
#    rename) even though the cell itself now reads "ProductsService", so
#    the hyperlink is added first and the cell value is re-applied after.
# ---------------------------------------------------------------------

$null = $ws.Hyperlinks.Add($ws.Range("D2"), "https://www.draw.io/", "G1J2xWywWx4tBD-WwT4hF07mk5S2uaKEWL")
$null = $ws.Hyperlinks.Add($ws.Range("D3"), "https://www.draw.io/", "G1oAa6kFLvZNQG8Zt8Qe23rxLx0LNsim0B")
$null = $ws.Hyperlinks.Add($ws.Range("D4"), "https://www.draw.io/", "G1cn6NAIUYOxWUz6hHUKNy_y4dG_9MVk_E", "", "InventoriesService")
$ws.Range("D4").Value = "ProductsService"

# Re-apply the built-in Hyperlink style explicitly so D2:D4 land on the
# same cell style as C2 (avoids a redundant near-duplicate style entry).
$ws.Range("D2:D4").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 4. Header row styling: C1:F1 bold on a themed green fill, with a thin
#    bottom border on C1:E1 only.
# ---------------------------------------------------------------------

$header = $ws.Range("C1:F1")
$header.Font.Bold = $true
$header.Interior.Color = 9555625

$headerBorder = $ws.Range("C1:E1")
$headerBorder.Borders.Item(9).LineStyle = 1
$headerBorder.Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------
# 5. Final selection, matching the saved cursor position in the diff.
# ---------------------------------------------------------------------

$null = $ws.Range("F16").Select()
